$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (style) of column H (header + data rows) onto the new
# columns I and J so the new cells match the existing look (header style s=1,
# plain data cells with no explicit style).
$ws.Range("H1:H31").Copy()
$ws.Range("I1:I31").PasteSpecial(-4122)
$ws.Range("H1:H31").Copy()
$ws.Range("J1:J31").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header row 1: new columns I ("I0") and J ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-31: new numeric values for columns I and J
$data = @{
    2  = @(7, 8)
    3  = @(8, 9)
    4  = @(7, 7)
    5  = @(7, 8)
    6  = @(7, 7)
    7  = @(6, 8)
    8  = @(8, 9)
    9  = @(8, 9)
    10 = @(7, 8)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(5, 8)
    14 = @(9, 9)
    15 = @(1, 6)
    16 = @(6, 7)
    17 = @(7, 8)
    18 = @(9, 9)
    19 = @(10, 10)
    20 = @(6, 8)
    21 = @(6, 8)
    22 = @(9, 9)
    23 = @(3, 6)
    24 = @(7, 8)
    25 = @(9, 9)
    26 = @(5, 7)
    27 = @(8, 8)
    28 = @(4, 6)
    29 = @(8, 8)
    30 = @(4, 6)
    31 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
